$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 460, shifting the
# existing rows 460-473 down to 462-475 (matching the target dimension
# A1:R475).
$ws.Rows.Item(460).Insert()
$ws.Rows.Item(460).Insert()

# Fill the two newly inserted rows with the new weekly price entries.
$ws.Range("A460").Value = 7
$ws.Range("B460").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C460").Value = "Ñuble"
$ws.Range("D460").Value2 = 45075
$ws.Range("E460").Value = 16
$ws.Range("F460").Value = 100112002
$ws.Range("G460").Value = "Pimiento"
$ws.Range("H460").Value = "Zafiro rojo"
$ws.Range("I460").Value = "Primera"
$ws.Range("J460").Value = 60
$ws.Range("K460").Value = 20000
$ws.Range("L460").Value = 20000
$ws.Range("M460").Value = 20000
$ws.Range("N460").Value = "$/caja 15 kilos"
$ws.Range("O460").Value = "Región de Arica y Parinacota"
$ws.Range("P460").Value = 1333
$ws.Range("Q460").Value = 15
$ws.Range("R460").Value = "Hortaliza"

$ws.Range("A461").Value = 7
$ws.Range("B461").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C461").Value = "Ñuble"
$ws.Range("D461").Value2 = 45075
$ws.Range("E461").Value = 16
$ws.Range("F461").Value = 100112002
$ws.Range("G461").Value = "Pimiento"
$ws.Range("H461").Value = "Zafiro verde"
$ws.Range("I461").Value = "Primera"
$ws.Range("J461").Value = 60
$ws.Range("K461").Value = 15000
$ws.Range("L461").Value = 15000
$ws.Range("M461").Value = 15000
$ws.Range("N461").Value = "$/caja 15 kilos"
$ws.Range("O461").Value = "Región de Arica y Parinacota"
$ws.Range("P461").Value = 1000
$ws.Range("Q461").Value = 15
$ws.Range("R461").Value = "Hortaliza"
